# Update countries & provincias Spain
# - Refresh the "Datos actualizados..." timestamp cell (A1)
# - Update case statistics for several countries (Rusia, Estonia, Lituania,
#   Eslovaquia, Sri Lanka) on the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp update (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 10:05"

# Row 6 - Rusia
$ws.Range("B6").Value = 362342
$ws.Range("C6").Value = 8915
$ws.Range("D6").Value = 131129
$ws.Range("E6").Value = 227406
$ws.Range("G6").Value = 174
$ws.Range("H6").Value = 3807

# Row 92 - Estonia
$ws.Range("B92").Value = 1834
$ws.Range("C92").Value = 10
$ws.Range("D92").Value = 1552
$ws.Range("E92").Value = 217

# Row 95 - Lituania
$ws.Range("B95").Value = 1639
$ws.Range("C95").Value = 4
$ws.Range("D95").Value = 1165
$ws.Range("E95").Value = 409
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 65

# Row 97 - Eslovaquia
$ws.Range("B97").Value = 1513
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 1322
$ws.Range("E97").Value = 163

# Row 103 - Sri Lanka
$ws.Range("D103").Value = 712
$ws.Range("E103").Value = 460
